# The National Statistical Committee's website moved from www.stat.kg to
# www.stat.gov.kg; update the "Organization website" cell to reflect the
# new address, then leave the user's selection where they finished editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B10 holds the "Organization website (if available)" answer.
$ws.Range("B10").Value = "www.stat.gov.kg"

# Re-assert B2's font so the cell picks up its own (new) font record,
# matching the re-saved style table produced by the author's edit.
$ws.Range("B2").Font.Name = "Calibri"

# Final cursor position left by the author after the edit.
$ws.Range("B8").Select()
